$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 50000  # H6: 10028.6 -> 50000
$ws.Cells.Item(6, 9).Value = 50000  # I6: 10028.6 -> 50000
$ws.Cells.Item(6, 11).Value = 150000  # K6: 30085.8 -> 150000
$ws.Cells.Item(6, 13).Value = -149888  # M6: -29973.8 -> -149888
$ws.Cells.Item(43, 8).Value = 9259951  # H43: 9259934 -> 9259951
$ws.Cells.Item(43, 9).Value = 475  # I43: 450.33334 -> 475
$ws.Cells.Item(43, 10).Value = 13889689  # J43: 18519418 -> 13889689
$ws.Cells.Item(43, 11).Value = 475  # K43: 450.33334 -> 475
$ws.Cells.Item(43, 12).Value = 13889689  # L43: 18519418 -> 13889689
$ws.Cells.Item(43, 13).Value = -406  # M43: -381.33334 -> -406
$ws.Cells.Item(43, 14).Value = -13889827  # N43: -18519556 -> -13889827
$ws.Cells.Item(64, 8).Value = 3553.6667  # H64: 3564.756 -> 3553.6667
$ws.Cells.Item(64, 9).Value = 3556.8948  # I64: 3582.3333 -> 3556.8948
$ws.Cells.Item(64, 11).Value = 3556.8948  # K64: 3582.3333 -> 3556.8948
$ws.Cells.Item(64, 13).Value = -3308.8948  # M64: -3334.3333 -> -3308.8948
$ws.Cells.Item(67, 8).Value = 3553.6667  # H67: 3564.756 -> 3553.6667
$ws.Cells.Item(67, 9).Value = 3556.8948  # I67: 3582.3333 -> 3556.8948
$ws.Cells.Item(67, 11).Value = 3556.8948  # K67: 3582.3333 -> 3556.8948
$ws.Cells.Item(67, 13).Value = -2698.8948  # M67: -2724.3333 -> -2698.8948
$ws.Cells.Item(74, 8).Value = 3483.3333  # H74: 3600 -> 3483.3333
$ws.Cells.Item(74, 9).Value = 3483.3333  # I74: 3600 -> 3483.3333
$ws.Cells.Item(74, 11).Value = 3483.3333  # K74: 3600 -> 3483.3333
$ws.Cells.Item(74, 13).Value = -2547.3333  # M74: -2664 -> -2547.3333
$ws.Cells.Item(77, 8).Value = 3483.3333  # H77: 3600 -> 3483.3333
$ws.Cells.Item(77, 9).Value = 3483.3333  # I77: 3600 -> 3483.3333
$ws.Cells.Item(77, 11).Value = 17416.6665  # K77: 18000 -> 17416.6665
$ws.Cells.Item(77, 13).Value = -12736.6665  # M77: -13320 -> -12736.6665
$ws.Cells.Item(112, 8).Value = 2570.7407  # H112: 2638.4614 -> 2570.7407
$ws.Cells.Item(112, 9).Value = 1099.5  # I112: 1100 -> 1099.5
$ws.Cells.Item(112, 10).Value = 2688.44  # J112: 2700 -> 2688.44
$ws.Cells.Item(112, 11).Value = 3298.5  # K112: 3300 -> 3298.5
$ws.Cells.Item(112, 12).Value = 8065.32  # L112: 8100 -> 8065.32
$ws.Cells.Item(112, 13).Value = -2190.5  # M112: -2192 -> -2190.5
$ws.Cells.Item(112, 14).Value = -10281.32  # N112: -10316 -> -10281.32
$ws.Cells.Item(116, 8).Value = 3575.1875  # H116: 3842.7144 -> 3575.1875
$ws.Cells.Item(116, 9).Value = 2819.9  # I116: 3099.25 -> 2819.9
$ws.Cells.Item(116, 11).Value = 2819.9  # K116: 3099.25 -> 2819.9
$ws.Cells.Item(116, 13).Value = 622.0999999999999  # M116: 342.75 -> 622.0999999999999
$ws.Cells.Item(137, 8).Value = 2146.3333  # H137: 1935.9688 -> 2146.3333
$ws.Cells.Item(137, 9).Value = 1735.6111  # I137: 1612.15 -> 1735.6111
$ws.Cells.Item(137, 10).Value = 2967.7778  # J137: 2475.6667 -> 2967.7778
$ws.Cells.Item(137, 11).Value = 5206.8333  # K137: 4836.450000000001 -> 5206.8333
$ws.Cells.Item(137, 12).Value = 8903.3334  # L137: 7427.000100000001 -> 8903.3334
$ws.Cells.Item(137, 13).Value = -2656.8333  # M137: -2286.450000000001 -> -2656.8333
$ws.Cells.Item(137, 14).Value = -14003.3334  # N137: -12527.0001 -> -14003.3334

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 300  # H5: 189 -> 300
$ws.Cells.Item(5, 9).Value = 300  # I5: 181.66667 -> 300
$ws.Cells.Item(5, 10).Value = 0  # J5: 200 -> 0
$ws.Cells.Item(5, 11).Value = 300  # K5: 181.66667 -> 300
$ws.Cells.Item(5, 12).Value = 0  # L5: 200 -> 0
$ws.Cells.Item(5, 13).Value = -188  # M5: -69.66667000000001 -> -188
$ws.Cells.Item(5, 14).ClearContents()  # N5 was -424
$ws.Cells.Item(32, 8).Value = 3116.6538  # H32: 2677.1638 -> 3116.6538
$ws.Cells.Item(32, 9).Value = 3303.5908  # I32: 2737.0186 -> 3303.5908
$ws.Cells.Item(32, 10).Value = 2088.5  # J32: 2215.4285 -> 2088.5
$ws.Cells.Item(32, 11).Value = 3303.5908  # K32: 2737.0186 -> 3303.5908
$ws.Cells.Item(32, 12).Value = 2088.5  # L32: 2215.4285 -> 2088.5
$ws.Cells.Item(32, 13).Value = -3016.5908  # M32: -2450.0186 -> -3016.5908
$ws.Cells.Item(32, 14).Value = -2662.5  # N32: -2789.4285 -> -2662.5
$ws.Cells.Item(88, 8).Value = 2684.5557  # H88: 2460 -> 2684.5557
$ws.Cells.Item(88, 9).Value = 0  # I88: 1793 -> 0
$ws.Cells.Item(88, 10).Value = 2684.5557  # J88: 2641.9092 -> 2684.5557
$ws.Cells.Item(88, 11).Value = 0  # K88: 1793 -> 0
$ws.Cells.Item(88, 12).Value = 2684.5557  # L88: 2641.9092 -> 2684.5557
$ws.Cells.Item(88, 13).ClearContents()  # M88 was -1387
$ws.Cells.Item(88, 14).Value = -3496.5557  # N88: -3453.9092 -> -3496.5557
$ws.Cells.Item(91, 8).Value = 2684.5557  # H91: 2460 -> 2684.5557
$ws.Cells.Item(91, 9).Value = 0  # I91: 1793 -> 0
$ws.Cells.Item(91, 10).Value = 2684.5557  # J91: 2641.9092 -> 2684.5557
$ws.Cells.Item(91, 11).Value = 0  # K91: 1793 -> 0
$ws.Cells.Item(91, 12).Value = 2684.5557  # L91: 2641.9092 -> 2684.5557
$ws.Cells.Item(91, 13).ClearContents()  # M91 was -389
$ws.Cells.Item(91, 14).Value = -5492.5557  # N91: -5449.9092 -> -5492.5557
$ws.Cells.Item(102, 8).Value = 83383340  # H102: 15161801 -> 83383340
$ws.Cells.Item(102, 9).Value = 83383340  # I102: 18530646 -> 83383340
$ws.Cells.Item(102, 10).Value = 0  # J102: 2000 -> 0
$ws.Cells.Item(102, 11).Value = 83383340  # K102: 18530646 -> 83383340
$ws.Cells.Item(102, 12).Value = 0  # L102: 2000 -> 0
$ws.Cells.Item(102, 13).Value = -83381718  # M102: -18529024 -> -83381718
$ws.Cells.Item(102, 14).ClearContents()  # N102 was -5244
$ws.Cells.Item(132, 8).Value = 3548.3462  # H132: 3739.88 -> 3548.3462
$ws.Cells.Item(132, 9).Value = 3433.2632  # I132: 3692.889 -> 3433.2632
$ws.Cells.Item(132, 11).Value = 10299.7896  # K132: 11078.667 -> 10299.7896
$ws.Cells.Item(132, 13).Value = -7769.7896  # M132: -8548.667000000001 -> -7769.7896

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 300  # H4: 189 -> 300
$ws.Cells.Item(4, 9).Value = 300  # I4: 181.66667 -> 300
$ws.Cells.Item(4, 10).Value = 0  # J4: 200 -> 0
$ws.Cells.Item(4, 11).Value = 300  # K4: 181.66667 -> 300
$ws.Cells.Item(4, 12).Value = 0  # L4: 200 -> 0
$ws.Cells.Item(4, 13).Value = -185  # M4: -66.66667000000001 -> -185
$ws.Cells.Item(4, 14).ClearContents()  # N4 was -430
$ws.Cells.Item(22, 8).Value = 371.33334  # H22: 271.25 -> 371.33334
$ws.Cells.Item(22, 9).Value = 345.6  # I22: 238.57143 -> 345.6
$ws.Cells.Item(22, 11).Value = 345.6  # K22: 238.57143 -> 345.6
$ws.Cells.Item(22, 13).Value = -172.6  # M22: -65.57142999999999 -> -172.6
$ws.Cells.Item(64, 8).Value = 231.61539  # H64: 320.83334 -> 231.61539
$ws.Cells.Item(64, 9).Value = 276  # I64: 320 -> 276
$ws.Cells.Item(64, 10).Value = 131.75  # J64: 321.66666 -> 131.75
$ws.Cells.Item(64, 11).Value = 276  # K64: 320 -> 276
$ws.Cells.Item(64, 12).Value = 131.75  # L64: 321.66666 -> 131.75
$ws.Cells.Item(64, 13).Value = -51  # M64: -95 -> -51
$ws.Cells.Item(64, 14).Value = -581.75  # N64: -771.66666 -> -581.75
$ws.Cells.Item(67, 8).Value = 231.61539  # H67: 320.83334 -> 231.61539
$ws.Cells.Item(67, 9).Value = 276  # I67: 320 -> 276
$ws.Cells.Item(67, 10).Value = 131.75  # J67: 321.66666 -> 131.75
$ws.Cells.Item(67, 11).Value = 276  # K67: 320 -> 276
$ws.Cells.Item(67, 12).Value = 131.75  # L67: 321.66666 -> 131.75
$ws.Cells.Item(67, 13).Value = 504  # M67: 460 -> 504
$ws.Cells.Item(67, 14).Value = -1691.75  # N67: -1881.66666 -> -1691.75
$ws.Cells.Item(86, 8).Value = 2932.7437  # H86: 3104.6 -> 2932.7437
$ws.Cells.Item(86, 9).Value = 3222.0833  # I86: 3505.4285 -> 3222.0833
$ws.Cells.Item(86, 10).Value = 2469.8  # J86: 2503.3572 -> 2469.8
$ws.Cells.Item(86, 11).Value = 3222.0833  # K86: 3505.4285 -> 3222.0833
$ws.Cells.Item(86, 12).Value = 2469.8  # L86: 2503.3572 -> 2469.8
$ws.Cells.Item(86, 13).Value = -2099.0833  # M86: -2382.4285 -> -2099.0833
$ws.Cells.Item(86, 14).Value = -4715.8  # N86: -4749.3572 -> -4715.8
$ws.Cells.Item(89, 8).Value = 2932.7437  # H89: 3104.6 -> 2932.7437
$ws.Cells.Item(89, 9).Value = 3222.0833  # I89: 3505.4285 -> 3222.0833
$ws.Cells.Item(89, 10).Value = 2469.8  # J89: 2503.3572 -> 2469.8
$ws.Cells.Item(89, 11).Value = 16110.4165  # K89: 17527.1425 -> 16110.4165
$ws.Cells.Item(89, 12).Value = 12349  # L89: 12516.786 -> 12349
$ws.Cells.Item(89, 13).Value = -10494.4165  # M89: -11911.1425 -> -10494.4165
$ws.Cells.Item(89, 14).Value = -23581  # N89: -23748.786 -> -23581
$ws.Cells.Item(134, 8).Value = 7312.5884  # H134: 8161.067 -> 7312.5884
$ws.Cells.Item(134, 9).Value = 1332.3846  # I134: 1402.091 -> 1332.3846
$ws.Cells.Item(134, 11).Value = 3997.1538  # K134: 4206.272999999999 -> 3997.1538
$ws.Cells.Item(134, 13).Value = -1462.1538  # M134: -1671.272999999999 -> -1462.1538

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 404  # H22: 397.5 -> 404
$ws.Cells.Item(22, 10).Value = 417.77777  # J22: 406 -> 417.77777
$ws.Cells.Item(22, 12).Value = 417.77777  # L22: 406 -> 417.77777
$ws.Cells.Item(22, 14).Value = -1117.77777  # N22: -1106 -> -1117.77777
$ws.Cells.Item(31, 8).Value = 1169.4928  # H31: 1199.7705 -> 1169.4928
$ws.Cells.Item(31, 9).Value = 787.7826  # I31: 816.6905 -> 787.7826
$ws.Cells.Item(31, 10).Value = 1932.9131  # J31: 2046.579 -> 1932.9131
$ws.Cells.Item(31, 11).Value = 787.7826  # K31: 816.6905 -> 787.7826
$ws.Cells.Item(31, 12).Value = 1932.9131  # L31: 2046.579 -> 1932.9131
$ws.Cells.Item(31, 13).Value = -492.7826  # M31: -521.6905 -> -492.7826
$ws.Cells.Item(31, 14).Value = -2522.9131  # N31: -2636.579 -> -2522.9131
$ws.Cells.Item(34, 8).Value = 1169.4928  # H34: 1199.7705 -> 1169.4928
$ws.Cells.Item(34, 9).Value = 787.7826  # I34: 816.6905 -> 787.7826
$ws.Cells.Item(34, 10).Value = 1932.9131  # J34: 2046.579 -> 1932.9131
$ws.Cells.Item(34, 11).Value = 787.7826  # K34: 816.6905 -> 787.7826
$ws.Cells.Item(34, 12).Value = 1932.9131  # L34: 2046.579 -> 1932.9131
$ws.Cells.Item(34, 13).Value = -585.7826  # M34: -614.6905 -> -585.7826
$ws.Cells.Item(34, 14).Value = -2336.9131  # N34: -2450.579 -> -2336.9131
$ws.Cells.Item(62, 8).Value = 5717243  # H62: 5266061 -> 5717243
$ws.Cells.Item(62, 9).Value = 3006.6333  # I62: 2980.6128 -> 3006.6333
$ws.Cells.Item(62, 10).Value = 40002664  # J62: 28573988 -> 40002664
$ws.Cells.Item(62, 11).Value = 3006.6333  # K62: 2980.6128 -> 3006.6333
$ws.Cells.Item(62, 12).Value = 40002664  # L62: 28573988 -> 40002664
$ws.Cells.Item(62, 13).Value = -2382.6333  # M62: -2356.6128 -> -2382.6333
$ws.Cells.Item(62, 14).Value = -40003912  # N62: -28575236 -> -40003912
$ws.Cells.Item(65, 8).Value = 5717243  # H65: 5266061 -> 5717243
$ws.Cells.Item(65, 9).Value = 3006.6333  # I65: 2980.6128 -> 3006.6333
$ws.Cells.Item(65, 10).Value = 40002664  # J65: 28573988 -> 40002664
$ws.Cells.Item(65, 11).Value = 15033.1665  # K65: 14903.064 -> 15033.1665
$ws.Cells.Item(65, 12).Value = 200013320  # L65: 142869940 -> 200013320
$ws.Cells.Item(65, 13).Value = -11913.1665  # M65: -11783.064 -> -11913.1665
$ws.Cells.Item(65, 14).Value = -200019560  # N65: -142876180 -> -200019560
$ws.Cells.Item(99, 8).Value = 1756017.6  # H99: 1755989.6 -> 1756017.6
$ws.Cells.Item(99, 9).Value = 2393950  # I99: 2194527.5 -> 2393950
$ws.Cells.Item(99, 10).Value = 1703.5  # J99: 1838 -> 1703.5
$ws.Cells.Item(99, 11).Value = 2393950  # K99: 2194527.5 -> 2393950
$ws.Cells.Item(99, 12).Value = 1703.5  # L99: 1838 -> 1703.5
$ws.Cells.Item(99, 13).Value = -2392452  # M99: -2193029.5 -> -2392452
$ws.Cells.Item(99, 14).Value = -4699.5  # N99: -4834 -> -4699.5
$ws.Cells.Item(126, 8).Value = 1756017.6  # H126: 1755989.6 -> 1756017.6
$ws.Cells.Item(126, 9).Value = 2393950  # I126: 2194527.5 -> 2393950
$ws.Cells.Item(126, 10).Value = 1703.5  # J126: 1838 -> 1703.5
$ws.Cells.Item(126, 11).Value = 7181850  # K126: 6583582.5 -> 7181850
$ws.Cells.Item(126, 12).Value = 5110.5  # L126: 5514 -> 5110.5
$ws.Cells.Item(126, 13).Value = -7179380  # M126: -6581112.5 -> -7179380
$ws.Cells.Item(126, 14).Value = -10050.5  # N126: -10454 -> -10050.5
$ws.Cells.Item(132, 8).Value = 5072.069  # H132: 5228.25 -> 5072.069
$ws.Cells.Item(132, 9).Value = 5038.9644  # I132: 5199.7036 -> 5038.9644
$ws.Cells.Item(132, 11).Value = 15116.8932  # K132: 15599.1108 -> 15116.8932
$ws.Cells.Item(132, 13).Value = -12586.8932  # M132: -13069.1108 -> -12586.8932

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 508253.8  # H4: 710534.4 -> 508253.8
$ws.Cells.Item(4, 9).Value = 100037.664  # I4: 81880.82000000001 -> 100037.664
$ws.Cells.Item(4, 10).Value = 701619.4  # J4: 1478888.6 -> 701619.4
$ws.Cells.Item(4, 11).Value = 300112.992  # K4: 245642.46 -> 300112.992
$ws.Cells.Item(4, 12).Value = 2104858.2  # L4: 4436665.800000001 -> 2104858.2
$ws.Cells.Item(4, 13).Value = -300000.992  # M4: -245530.46 -> -300000.992
$ws.Cells.Item(4, 14).Value = -2105082.2  # N4: -4436889.800000001 -> -2105082.2
$ws.Cells.Item(36, 8).Value = 500  # H36: 0 -> 500
$ws.Cells.Item(36, 9).Value = 500  # I36: 0 -> 500
$ws.Cells.Item(36, 11).Value = 1500  # K36: 0 -> 1500
$ws.Cells.Item(36, 13).Value = -1331  # M36: None -> -1331
$ws.Cells.Item(68, 8).Value = 2170.7273  # H68: 2295.9412 -> 2170.7273
$ws.Cells.Item(68, 10).Value = 2170.7273  # J68: 2295.9412 -> 2170.7273
$ws.Cells.Item(68, 12).Value = 6512.1819  # L68: 6887.823600000001 -> 6512.1819
$ws.Cells.Item(68, 14).Value = -8134.1819  # N68: -8509.8236 -> -8134.1819
$ws.Cells.Item(71, 8).Value = 2170.7273  # H71: 2295.9412 -> 2170.7273
$ws.Cells.Item(71, 10).Value = 2170.7273  # J71: 2295.9412 -> 2170.7273
$ws.Cells.Item(71, 12).Value = 19536.5457  # L71: 20663.4708 -> 19536.5457
$ws.Cells.Item(71, 14).Value = -27648.5457  # N71: -28775.4708 -> -27648.5457
$ws.Cells.Item(98, 8).Value = 218.33333  # H98: 40 -> 218.33333
$ws.Cells.Item(98, 9).Value = 177.5  # I98: 0 -> 177.5
$ws.Cells.Item(98, 10).Value = 300  # J98: 40 -> 300
$ws.Cells.Item(98, 11).Value = 532.5  # K98: 0 -> 532.5
$ws.Cells.Item(98, 12).Value = 900  # L98: 120 -> 900
$ws.Cells.Item(98, 13).Value = 965.5  # M98: None -> 965.5
$ws.Cells.Item(98, 14).Value = -3896  # N98: -3116 -> -3896
$ws.Cells.Item(107, 8).Value = 7981.643  # H107: 10076.637 -> 7981.643
$ws.Cells.Item(107, 9).Value = 675.5  # I107: 1051 -> 675.5
$ws.Cells.Item(107, 11).Value = 2026.5  # K107: 3153 -> 2026.5
$ws.Cells.Item(107, 13).Value = -106.5  # M107: -1233 -> -106.5
$ws.Cells.Item(112, 8).Value = 10800.267  # H112: 47627936 -> 10800.267
$ws.Cells.Item(112, 9).Value = 2150.6667  # I112: 1780.8 -> 2150.6667
$ws.Cells.Item(112, 10).Value = 16566.666  # J112: 62511108 -> 16566.666
$ws.Cells.Item(112, 11).Value = 6452.000100000001  # K112: 5342.4 -> 6452.000100000001
$ws.Cells.Item(112, 12).Value = 49699.99800000001  # L112: 187533324 -> 49699.99800000001
$ws.Cells.Item(112, 13).Value = -5344.000100000001  # M112: -4234.4 -> -5344.000100000001
$ws.Cells.Item(112, 14).Value = -51915.99800000001  # N112: -187535540 -> -51915.99800000001
$ws.Cells.Item(113, 8).Value = 597.3714  # H113: 606.0294 -> 597.3714
$ws.Cells.Item(113, 9).Value = 498.36365  # I113: 517.9 -> 498.36365
$ws.Cells.Item(113, 11).Value = 1495.09095  # K113: 1553.7 -> 1495.09095
$ws.Cells.Item(113, 13).Value = 674.90905  # M113: 616.3000000000002 -> 674.90905
$ws.Cells.Item(122, 8).Value = 547  # H122: 489.76923 -> 547
$ws.Cells.Item(122, 9).Value = 402.66666  # I122: 367.4 -> 402.66666
$ws.Cells.Item(122, 10).Value = 1196.5  # J122: 897.6667 -> 1196.5
$ws.Cells.Item(122, 11).Value = 3623.99994  # K122: 3306.6 -> 3623.99994
$ws.Cells.Item(122, 12).Value = 10768.5  # L122: 8079.0003 -> 10768.5
$ws.Cells.Item(122, 13).Value = -1173.99994  # M122: -856.5999999999999 -> -1173.99994
$ws.Cells.Item(122, 14).Value = -15668.5  # N122: -12979.0003 -> -15668.5
$ws.Cells.Item(137, 8).Value = 16973.223  # H137: 15368.9 -> 16973.223
$ws.Cells.Item(137, 9).Value = 3165  # I137: 2718 -> 3165
$ws.Cells.Item(137, 11).Value = 9495  # K137: 8154 -> 9495
$ws.Cells.Item(137, 13).Value = -4395  # M137: -3054 -> -4395
$ws.Cells.Item(140, 8).Value = 32856.5  # H140: 27198.045 -> 32856.5
$ws.Cells.Item(140, 9).Value = 38834.934  # I140: 37600.902 -> 38834.934
$ws.Cells.Item(140, 10).Value = 2964.3333  # J140: 2391.2307 -> 2964.3333
$ws.Cells.Item(140, 11).Value = 116504.802  # K140: 112802.706 -> 116504.802
$ws.Cells.Item(140, 12).Value = 8892.999899999999  # L140: 7173.6921 -> 8892.999899999999
$ws.Cells.Item(140, 13).Value = -111324.802  # M140: -107622.706 -> -111324.802
$ws.Cells.Item(140, 14).Value = -19252.9999  # N140: -17533.6921 -> -19252.9999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(74, 8).Value = 66200  # H74: 63239.668 -> 66200
$ws.Cells.Item(74, 9).Value = 0  # I74: 50119 -> 0
$ws.Cells.Item(74, 10).Value = 66200  # J74: 69800 -> 66200
$ws.Cells.Item(74, 11).Value = 0  # K74: 50119 -> 0
$ws.Cells.Item(74, 12).Value = 66200  # L74: 69800 -> 66200
$ws.Cells.Item(74, 13).ClearContents()  # M74 was -49183
$ws.Cells.Item(74, 14).Value = -68072  # N74: -71672 -> -68072
$ws.Cells.Item(77, 8).Value = 66200  # H77: 63239.668 -> 66200
$ws.Cells.Item(77, 9).Value = 0  # I77: 50119 -> 0
$ws.Cells.Item(77, 10).Value = 66200  # J77: 69800 -> 66200
$ws.Cells.Item(77, 11).Value = 0  # K77: 150357 -> 0
$ws.Cells.Item(77, 12).Value = 198600  # L77: 209400 -> 198600
$ws.Cells.Item(77, 13).ClearContents()  # M77 was -145677
$ws.Cells.Item(77, 14).Value = -207960  # N77: -218760 -> -207960

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2600.6667  # H40: 3101.3333 -> 2600.6667
$ws.Cells.Item(40, 9).Value = 2426  # I40: 2952 -> 2426
$ws.Cells.Item(40, 10).Value = 2950  # J40: 3400 -> 2950
$ws.Cells.Item(40, 11).Value = 2426  # K40: 2952 -> 2426
$ws.Cells.Item(40, 12).Value = 2950  # L40: 3400 -> 2950
$ws.Cells.Item(40, 13).Value = -2290  # M40: -2816 -> -2290
$ws.Cells.Item(40, 14).Value = -3222  # N40: -3672 -> -3222
$ws.Cells.Item(132, 8).Value = 31412.324  # H132: 34352.355 -> 31412.324
$ws.Cells.Item(132, 9).Value = 1632.2106  # I132: 1642.9474 -> 1632.2106
$ws.Cells.Item(132, 10).Value = 69133.8  # J132: 86142.25 -> 69133.8
$ws.Cells.Item(132, 11).Value = 4896.6318  # K132: 4928.8422 -> 4896.6318
$ws.Cells.Item(132, 12).Value = 207401.4  # L132: 258426.75 -> 207401.4
$ws.Cells.Item(132, 13).Value = -2366.6318  # M132: -2398.8422 -> -2366.6318
$ws.Cells.Item(132, 14).Value = -212461.4  # N132: -263486.75 -> -212461.4
$ws.Cells.Item(136, 8).Value = 6610.35  # H136: 7212.0557 -> 6610.35
$ws.Cells.Item(136, 9).Value = 7260.647  # I136: 8069.4 -> 7260.647
$ws.Cells.Item(136, 11).Value = 21781.941  # K136: 24208.2 -> 21781.941
$ws.Cells.Item(136, 13).Value = -19231.941  # M136: -21658.2 -> -19231.941

Write-Host "Applied all changes"